$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each target cell is forced to Text format before assignment so that
# numeric-looking strings (e.g. "41.70", "1.001") are stored as literal
# text instead of being auto-coerced to a Double by the input parser.
# The NumberFormat is then reset back to the default "Normal" style so
# the cell's style index is unchanged from the original file.
$cellValues = [ordered]@{
    'D2' = '28.061.27'
    'E2' = '  -1.15%  '
    'D3' = '1.793.69'
    'E3' = '  -0.43%  '
    'E4' = '  +0.00%  '
    'D5' = '316.93'
    'E5' = '  +0.65%  '
    'E6' = '  +0.07%  '
    'D7' = '0.5390'
    'E7' = '  -2.14%  '
    'D8' = '0.3768'
    'E8' = '  -2.26%  '
    'D9' = '0.07451'
    'E9' = '  -1.96%  '
    'D10' = '41.70'
    'E10' = '  -1.98%  '
    'E11' = '  -3.06%  '
    'D12' = '1.001'
    'E12' = '  +0.02%  '
    'D13' = '20.56'
    'E13' = '  -3.03%  '
    'D14' = '6.095'
    'E14' = '  -1.44%  '
    'D15' = '1.786.62'
    'E15' = '  -1.29%  '
    'D16' = '7.209'
    'E16' = '  -2.87%  '
    'D17' = '89.02'
    'E17' = '  -3.27%  '
    'E18' = '  -1.56%  '
    'D19' = '0.06466'
    'E19' = '  +0.36%  '
    'E20' = '  +0.05%  '
    'D21' = '17.27'
    'E21' = '  -0.54%  '
    'D22' = '5.891'
    'E22' = '  -1.38%  '
    'D23' = '28.090.01'
    'E23' = '  -1.20%  '
    'D24' = '11.16'
    'E24' = '  -2.49%  '
    'D25' = '2.094'
    'E25' = '  -1.84%  '
    'D26' = '154.80'
    'E26' = '  -2.53%  '
    'E27' = '  -2.31%  '
    'D28' = '1.991.20'
    'E28' = '  -1.26%  '
    'E29' = '  -5.48%  '
    'D30' = '120.56'
    'E30' = '  -2.76%  '
    'D31' = '1.115'
    'E31' = '  -0.59%  '
    'D32' = '0.1056'
    'E32' = '  +3.04%  '
    'E33' = '  -0.94%  '
    'D34' = '5.545'
    'E34' = '  -3.98%  '
    'D35' = '0.06525'
    'E35' = '  +1.91%  '
    'D36' = '0.2256'
    'E36' = '  -2.44%  '
    'D37' = '0.02279'
    'E37' = '  -2.06%  '
    'D38' = '5.012'
    'E38' = '  -2.85%  '
    'D39' = '8.437'
    'E39' = '  -3.85%  '
    'B40' = 'WEMIXTOKEN'
    'C40' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D40' = '1.446'
    'E40' = '  +4.49%  '
    'B41' = 'TheSandbox'
    'C41' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D41' = '0.6158'
    'E41' = '  -4.13%  '
    'E42' = '  -5.06%  '
    'D43' = '1.170'
    'E43' = '  +0.91%  '
    'D44' = '0.9998'
    'E44' = '  +0.02%  '
    'D45' = '13.35'
    'E45' = '  -1.79%  '
    'D46' = '3.672'
    'E46' = '  -0.22%  '
    'D47' = '0.5768'
    'E47' = '  -3.57%  '
    'D48' = '125.94'
    'E48' = '  -0.59%  '
    'D49' = '1.186'
    'E49' = '  +3.00%  '
    'D50' = '1.916'
    'E50' = '  -3.48%  '
    'D51' = '0.06812'
    'E51' = '  -1.27%  '
}

foreach ($ref in $cellValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cellValues[$ref]
    $cell.Style = "Normal"
}
